$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-05 -> 2023-10-06) for every data row (2 through 211).
for ($row = 2; $row -le 211; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
